$d = $word.ActiveDocument

$replacements = @(
    @{old="147÷4="; new="839÷9="},
    @{old="360÷5="; new="963÷4="},
    @{old="737÷8="; new="907÷6="},
    @{old="182÷6="; new="177÷5="},
    @{old="969÷3="; new="302÷4="},
    @{old="496÷5="; new="297÷5="},
    @{old="296÷9="; new="668÷9="},
    @{old="827÷3="; new="250÷8="},
    @{old="146÷7="; new="422÷2="},
    @{old="278÷9="; new="564÷3="},
    @{old="265÷6="; new="927÷5="},
    @{old="633÷6="; new="294÷6="},
    @{old="684÷8="; new="190÷7="},
    @{old="544÷8="; new="557÷7="},
    @{old="514÷3="; new="509÷4="},
    @{old="374÷6="; new="334÷2="},
    @{old="474÷2="; new="642÷8="},
    @{old="641÷2="; new="345÷9="},
    @{old="727÷3="; new="562÷3="},
    @{old="693÷7="; new="520÷8="},
    @{old="623÷6="; new="813÷6="},
    @{old="678÷2="; new="119÷2="},
    @{old="334÷6="; new="694÷7="},
    @{old="258÷5="; new="198÷4="},
    @{old="983÷2="; new="168÷3="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
